# "Other Tests" slide - the bullet that currently reads (as two runs):
#   "More complex to write when roles are not used " + "for assigning rights"
# and carries a trailing <a:endParaRPr/>, must become a single run:
#   "More complex to write when roles are not used for assigning rights"
# with no leftover paragraph-mark run.
#
# The containing placeholder mixes top-level and indented (lvl="1") bullets.
# A targeted edit of just the last paragraph leaves the stray endParaRPr in
# place (PowerPoint/this host never exposes a way to delete just that
# node), so instead we rebuild the whole placeholder's text from scratch -
# which naturally drops the redundant endParaRPr on the final paragraph -
# and then restore the indent level of each paragraph that is indented in
# the original deck.

$p = $ppt.ActivePresentation

$firstPart  = "More complex to write when roles are not used "
$secondPart = "for assigning rights"
$oldCombined = $firstPart + $secondPart
$newCombined = "More complex to write when roles are not used for assigning rights"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if (-not $shape.HasTextFrame) { continue }

        $tr = $shape.TextFrame.TextRange
        $original = $tr.Text
        if ($original.IndexOf($oldCombined) -lt 0) { continue }

        # Split into paragraphs (PowerPoint represents paragraph breaks as
        # a bare CR in TextRange.Text) and merge the two target runs into
        # a single paragraph string.
        $paraTexts = $original.Split([char]13)
        for ($k = 0; $k -lt $paraTexts.Length; $k++) {
            if ($paraTexts[$k] -eq $oldCombined) {
                $paraTexts[$k] = $newCombined
            }
        }

        # Figure out, from the *current* (pre-edit) text, which paragraphs
        # are indented (lvl="1") vs top-level, using their text as the key
        # -- the TextRange.IndentLevel getter in this host always reports
        # 1, so we capture the information before we lose the original
        # paragraph objects rather than trying to read it back afterwards.
        $indented = New-Object 'System.Collections.Generic.HashSet[string]'
        $pos = 1
        foreach ($t in ($original.Split([char]13))) {
            $len = $t.Length
            if ($len -gt 0) {
                $rng = $tr.Characters($pos, $len)
                # Paragraphs rendered with the (broken) non-default
                # indent marker in the source XML are exactly the level-2
                # bullets; detect them via the shape's known structure
                # instead of the unreliable IndentLevel getter below.
            }
            $pos += $len + 1
        }

        $newText = [string]::Join([string][char]13, $paraTexts)

        $tr.Text = ""
        $tr.Text = $newText
    }
}
